# This workbook tracks daily fruit/vegetable price records, one row per
# quality-grade line item. The edit inserts 3 new records (a new pricing
# date, 44644) right above the existing "2021-10-14" (44483) block for
# "1a/2a/3a amarillo", pushing every subsequent row down by 3 and growing
# the used range from A1:T538 to A1:T541.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert 3 new blank rows just above row 457 --------------------------
$ws.Range("A457:A459").EntireRow.Insert()

# --- Populate the 3 new rows with the new price record ------------------
# Common (unchanged) values across the 3 new rows
$mercadoId = 2
$mercado   = "Comercializadora del Agro de Limarí"
$region    = "Coquimbo"
$fecha     = 44644
$codreg    = 4
$tipo      = "Fruta"
$prodId    = 100102
$producto  = "Cítricos"
$catId     = 100102003
$categoria = "Limón"
$variedad  = "Sin especificar"
$unidad    = "$/malla 16 kilos"
$origen    = "Provincia de Limarí"
$kgUnidad  = 16

# Row 457: 1a amarillo
$ws.Range("A457").Value2 = $mercadoId
$ws.Range("B457").Value2 = $mercado
$ws.Range("C457").Value2 = $region
$ws.Range("D457").Value2 = $fecha
$ws.Range("E457").Value2 = $codreg
$ws.Range("F457").Value2 = $tipo
$ws.Range("G457").Value2 = $prodId
$ws.Range("H457").Value2 = $producto
$ws.Range("I457").Value2 = $catId
$ws.Range("J457").Value2 = $categoria
$ws.Range("K457").Value2 = $variedad
$ws.Range("L457").Value2 = "1a amarillo"
$ws.Range("M457").Value2 = 750
$ws.Range("N457").Value2 = 17800
$ws.Range("O457").Value2 = 18000
$ws.Range("P457").Value2 = 17900
$ws.Range("Q457").Value2 = $unidad
$ws.Range("R457").Value2 = $origen
$ws.Range("S457").Value2 = 1119
$ws.Range("T457").Value2 = $kgUnidad

# Row 458: 2a amarillo
$ws.Range("A458").Value2 = $mercadoId
$ws.Range("B458").Value2 = $mercado
$ws.Range("C458").Value2 = $region
$ws.Range("D458").Value2 = $fecha
$ws.Range("E458").Value2 = $codreg
$ws.Range("F458").Value2 = $tipo
$ws.Range("G458").Value2 = $prodId
$ws.Range("H458").Value2 = $producto
$ws.Range("I458").Value2 = $catId
$ws.Range("J458").Value2 = $categoria
$ws.Range("K458").Value2 = $variedad
$ws.Range("L458").Value2 = "2a amarillo"
$ws.Range("M458").Value2 = 600
$ws.Range("N458").Value2 = 14800
$ws.Range("O458").Value2 = 15000
$ws.Range("P458").Value2 = 14900
$ws.Range("Q458").Value2 = $unidad
$ws.Range("R458").Value2 = $origen
$ws.Range("S458").Value2 = 931
$ws.Range("T458").Value2 = $kgUnidad

# Row 459: 3a amarillo
$ws.Range("A459").Value2 = $mercadoId
$ws.Range("B459").Value2 = $mercado
$ws.Range("C459").Value2 = $region
$ws.Range("D459").Value2 = $fecha
$ws.Range("E459").Value2 = $codreg
$ws.Range("F459").Value2 = $tipo
$ws.Range("G459").Value2 = $prodId
$ws.Range("H459").Value2 = $producto
$ws.Range("I459").Value2 = $catId
$ws.Range("J459").Value2 = $categoria
$ws.Range("K459").Value2 = $variedad
$ws.Range("L459").Value2 = "3a amarillo"
$ws.Range("M459").Value2 = 420
$ws.Range("N459").Value2 = 12800
$ws.Range("O459").Value2 = 13000
$ws.Range("P459").Value2 = 12900
$ws.Range("Q459").Value2 = $unidad
$ws.Range("R459").Value2 = $origen
$ws.Range("S459").Value2 = 806
$ws.Range("T459").Value2 = $kgUnidad

# Match the date-column number format used by the rest of column D so the
# new cells render as dates rather than raw serials.
$ws.Range("D457:D459").NumberFormat = $ws.Range("D460").NumberFormat()
